# Scheduled runner update: refresh cached market-board prices / leve profit
# figures across all eight Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below are the new snapshot pulled for each leve row; a handful of
# rows gain or lose their HQ-profit (N) / NQ-profit (M) cell entirely when
# the underlying recipe no longer prices out an HQ (or NQ) variant.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 54.81818
$ws.Range("I9").Value = 61.333332
$ws.Range("J9").Value = 25.5
$ws.Range("K9").Value = 61.333332
$ws.Range("L9").Value = 25.5
$ws.Range("M9").Value = 107.666668
$ws.Range("N9").Value = -363.5

$ws.Range("H17").Value = 3379.75
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 3655.7
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 10967.1
$ws.Range("M17").Value = -5832
$ws.Range("N17").Value = -11303.1

$ws.Range("H115").Value = 1512.8125
$ws.Range("I115").Value = 959.2857
$ws.Range("K115").Value = 2877.8571
$ws.Range("M115").Value = -1310.8571

$ws.Range("H132").Value = 11793.667
$ws.Range("I132").Value = 12182.685
$ws.Range("J132").Value = 10315.4
$ws.Range("K132").Value = 36548.055
$ws.Range("L132").Value = 30946.2
$ws.Range("M132").Value = -34018.055
$ws.Range("N132").Value = -36006.2

$ws.Range("H137").Value = 2000
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3450
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 6426.2666
$ws.Range("I138").Value = 2732.6667
$ws.Range("K138").Value = 8198.000100000001
$ws.Range("M138").Value = -3058.000100000001

$ws.Range("H140").Value = 87390
$ws.Range("J140").Value = 87390
$ws.Range("L140").Value = 87390
$ws.Range("N140").Value = -97750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 450
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -674

$ws.Range("H61").Value = 5538.857
$ws.Range("I61").Value = 5349.5386
$ws.Range("K61").Value = 5349.5386
$ws.Range("M61").Value = -5137.5386

$ws.Range("H74").Value = 5000.8335
$ws.Range("I74").Value = 3801.2
$ws.Range("J74").Value = 10999
$ws.Range("K74").Value = 3801.2
$ws.Range("L74").Value = 10999
$ws.Range("M74").Value = -2927.2
$ws.Range("N74").Value = -12747

$ws.Range("H77").Value = 5000.8335
$ws.Range("I77").Value = 3801.2
$ws.Range("J77").Value = 10999
$ws.Range("K77").Value = 19006
$ws.Range("L77").Value = 54995
$ws.Range("M77").Value = -14638
$ws.Range("N77").Value = -63731

$ws.Range("H136").Value = 5538.857
$ws.Range("I136").Value = 5349.5386
$ws.Range("K136").Value = 16048.6158
$ws.Range("M136").Value = -13498.6158

$ws.Range("H137").Value = 65000
$ws.Range("J137").Value = 65000
$ws.Range("L137").Value = 65000
$ws.Range("N137").Value = -75200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -680

$ws.Range("H107").Value = 2832.5
$ws.Range("I107").Value = 1376.625
$ws.Range("K107").Value = 1376.625
$ws.Range("M107").Value = 543.375

$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -110119

$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3207.8572
$ws.Range("I16").Value = 2988.6667
$ws.Range("K16").Value = 2988.6667
$ws.Range("M16").Value = -2701.6667

$ws.Range("H22").Value = 4576
$ws.Range("I22").Value = 4224.5
$ws.Range("J22").Value = 4732.222
$ws.Range("K22").Value = 4224.5
$ws.Range("L22").Value = 4732.222
$ws.Range("M22").Value = -3874.5
$ws.Range("N22").Value = -5432.222

$ws.Range("H58").Value = 1937.75
$ws.Range("I58").Value = 1250
$ws.Range("K58").Value = 1250
$ws.Range("M58").Value = -1047

$ws.Range("H59").Value = 60491.8
$ws.Range("I59").Value = 19904
$ws.Range("K59").Value = 19904
$ws.Range("M59").Value = -18759

$ws.Range("H60").Value = 878.6
$ws.Range("I60").Value = 878.6
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 878.6
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -367.6
$ws.Range("N60").ClearContents()

$ws.Range("H113").Value = 3207.8572
$ws.Range("I113").Value = 2988.6667
$ws.Range("K113").Value = 2988.6667
$ws.Range("M113").Value = -818.6667000000002

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H132").Value = 3971.6667
$ws.Range("I132").Value = 3035.3572
$ws.Range("K132").Value = 9106.071599999999
$ws.Range("M132").Value = -6576.071599999999

$ws.Range("H134").Value = 2710.8667
$ws.Range("I134").Value = 2130.3076
$ws.Range("K134").Value = 6390.9228
$ws.Range("M134").Value = -3855.9228

$ws.Range("H136").Value = 1937.75
$ws.Range("I136").Value = 1250
$ws.Range("K136").Value = 3750
$ws.Range("M136").Value = -1200

$ws.Range("H141").Value = 48058.465
$ws.Range("J141").Value = 49348.355
$ws.Range("L141").Value = 49348.355
$ws.Range("N141").Value = -59708.355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1590
$ws.Range("I18").Value = 1036.2
$ws.Range("K18").Value = 3108.6
$ws.Range("M18").Value = -2939.6

$ws.Range("H60").Value = 603.9091
$ws.Range("I60").Value = 277.94736
$ws.Range("K60").Value = 833.84208
$ws.Range("M60").Value = -582.84208

$ws.Range("H122").Value = 889.6667
$ws.Range("I122").Value = 624
$ws.Range("J122").Value = 1421
$ws.Range("K122").Value = 5616
$ws.Range("L122").Value = 12789
$ws.Range("M122").Value = -3166
$ws.Range("N122").Value = -17689

$ws.Range("H124").Value = 5030.25
$ws.Range("I124").Value = 5000
$ws.Range("K124").Value = 15000
$ws.Range("M124").Value = -10090

$ws.Range("H131").Value = 1531
$ws.Range("J131").Value = 1397.2
$ws.Range("L131").Value = 4191.6
$ws.Range("N131").Value = -14271.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2199.5
$ws.Range("I80").Value = 1904
$ws.Range("J80").Value = 2495
$ws.Range("K80").Value = 1904
$ws.Range("L80").Value = 2495
$ws.Range("M80").Value = -906
$ws.Range("N80").Value = -4491

$ws.Range("H83").Value = 2199.5
$ws.Range("I83").Value = 1904
$ws.Range("J83").Value = 2495
$ws.Range("K83").Value = 9520
$ws.Range("L83").Value = 12475
$ws.Range("M83").Value = -4528
$ws.Range("N83").Value = -22459

$ws.Range("H93").Value = 56999.5
$ws.Range("J93").Value = 56999.5
$ws.Range("L93").Value = 56999.5
$ws.Range("N93").Value = -60743.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5025.8
$ws.Range("I7").Value = 2761.8
$ws.Range("K7").Value = 2761.8
$ws.Range("M7").Value = -2649.8

$ws.Range("H61").Value = 5394.4443
$ws.Range("I61").Value = 4550
$ws.Range("K61").Value = 4550
$ws.Range("M61").Value = -4348

$ws.Range("H113").Value = 5394.4443
$ws.Range("I113").Value = 4550
$ws.Range("K113").Value = 4550
$ws.Range("M113").Value = -2380

$ws.Range("H126").Value = 5025.8
$ws.Range("I126").Value = 2761.8
$ws.Range("K126").Value = 8285.400000000001
$ws.Range("M126").Value = -5815.400000000001

$ws.Range("H132").Value = 9422.477000000001
$ws.Range("I132").Value = 8205.429
$ws.Range("J132").Value = 11856.571
$ws.Range("K132").Value = 24616.287
$ws.Range("L132").Value = 35569.713
$ws.Range("M132").Value = -22086.287
$ws.Range("N132").Value = -40629.713

$ws.Range("H136").Value = 4699.75
$ws.Range("I136").Value = 4219.6
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 12658.8
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -10108.8
$ws.Range("N136").Value = -21600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 367887
$ws.Range("I4").Value = 471497.56
$ws.Range("K4").Value = 471497.56
$ws.Range("M4").Value = -471384.56

$ws.Range("H132").Value = 2599.9167
$ws.Range("I132").Value = 2599.9167
$ws.Range("K132").Value = 7799.750100000001
$ws.Range("M132").Value = -5269.750100000001

$ws.Range("H136").Value = 6183.136
$ws.Range("I136").Value = 5336.222
$ws.Range("K136").Value = 16008.666
$ws.Range("M136").Value = -13458.666
